$d = $word.ActiveDocument

# --- Paragraph 5 ("Component" + "s " + "Render other Components") ---
# Collapses the 3 runs into a single run reading "this.props", wrapped in
# proofErr spellStart/spellEnd markers (as Word's proofer would do for an
# unrecognised "word" like this.props).
$p5 = $d.Paragraphs.Item(5)
$r5 = $p5.Range
$x5 = $r5.WordOpenXML
$old5 = '<w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas"/><w:color w:val="00B0F0"/><w:sz w:val="48"/><w:szCs w:val="48"/></w:rPr><w:t>Components Render other Components</w:t></w:r>'
$new5 = '<w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas"/><w:color w:val="00B0F0"/><w:sz w:val="48"/><w:szCs w:val="48"/></w:rPr><w:t>this.props</w:t></w:r><w:proofErr w:type="spellEnd"/>'
$x5n = $x5.Replace($old5, $new5)
$r5.InsertXML($x5n)

# --- Paragraph 6 ("Components Interact") ---
# Retargets the run text to "this.props", wrapped in proofErr markers.
$p6 = $d.Paragraphs.Item(6)
$r6 = $p6.Range
$x6 = $r6.WordOpenXML
$old6 = '<w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas"/><w:color w:val="00B050"/><w:sz w:val="36"/><w:szCs w:val="36"/></w:rPr><w:t>Components Interact</w:t></w:r>'
$new6 = '<w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas"/><w:color w:val="00B050"/><w:sz w:val="36"/><w:szCs w:val="36"/></w:rPr><w:t>this.props</w:t></w:r><w:proofErr w:type="spellEnd"/>'
$x6n = $x6.Replace($old6, $new6)
$r6.InsertXML($x6n)

# --- Remove the 4 blank spacer paragraphs right after "this.props" ---
# (originally paragraphs 7-10, the empty Consolas/sz24 paragraphs); the
# 3 blank paragraphs that follow them are left untouched.
for ($i = 1; $i -le 4; $i++) {
    $blank = $d.Paragraphs.Item(7)
    $blank.Range.Delete()
}
